$d = $word.ActiveDocument

$pairs = @(
    @("904÷3=", "241÷6="),
    @("882÷2=", "185÷6="),
    @("425÷6=", "268÷9="),
    @("364÷3=", "953÷7="),
    @("411÷9=", "715÷5="),
    @("912÷9=", "749÷3="),
    @("289÷7=", "787÷6="),
    @("906÷7=", "481÷2="),
    @("698÷9=", "689÷5="),
    @("862÷8=", "531÷3="),
    @("104÷6=", "803÷8="),
    @("653÷3=", "607÷3="),
    @("335÷6=", "635÷8="),
    @("615÷7=", "407÷3="),
    @("618÷7=", "305÷9="),
    @("241÷4=", "650÷7="),
    @("502÷8=", "370÷3="),
    @("279÷7=", "319÷8="),
    @("101÷5=", "311÷6="),
    @("701÷8=", "978÷8="),
    @("703÷2=", "286÷6="),
    @("974÷7=", "642÷7="),
    @("956÷8=", "333÷4="),
    @("984÷5=", "145÷8="),
    @("206÷3=", "896÷5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
